$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: was Hedera -> becomes Dai
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.08%  "

# Row 36: was Dai -> becomes Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.865.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.268.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.266.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "685.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.799.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.088.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.277.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.878"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("E26").Value = "  -5.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "576.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.798.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -16.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.127"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.322"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0400"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
